{"js": "// Update the three-digit-by-one-digit division answers in the table.\n// Each \"old\" answer string is unique within the document, so a plain\n// text search + full-text replace of each matching range is safe and\n// unambiguous (one-to-one mapping, order does not matter).\nconst replacements = [\n  [\"766\u00f78=95, 6\", \"892\u00f75=178, 2\"],\n  [\"916\u00f76=152, 4\", \"956\u00f73=318, 2\"],\n  [\"375\u00f74=93, 3\", \"644\u00f72=322, 0\"],\n  [\"825\u00f75=165, 0\", \"545\u00f76=90, 5\"],\n  [\"936\u00f76=156, 0\", \"204\u00f74=51, 0\"],\n  [\"279\u00f77=39, 6\", \"250\u00f75=50, 0\"],\n  [\"851\u00f76=141, 5\", \"788\u00f78=98, 4\"],\n  [\"351\u00f74=87, 3\", \"781\u00f77=111, 4\"],\n  [\"152\u00f73=50, 2\", \"428\u00f73=142, 2\"],\n  [\"692\u00f76=115, 2\", \"672\u00f74=168, 0\"],\n  [\"432\u00f75=86, 2\", \"514\u00f79=57, 1\"],\n  [\"772\u00f72=386, 0\", \"837\u00f77=119, 4\"],\n  [\"791\u00f72=395, 1\", \"114\u00f78=14, 2\"],\n  [\"498\u00f74=124, 2\", \"457\u00f72=228, 1\"],\n  [\"786\u00f76=131, 0\", \"865\u00f74=216, 1\"],\n  [\"116\u00f79=12, 8\", \"515\u00f76=85, 5\"],\n  [\"583\u00f79=64, 7\", \"398\u00f73=132, 2\"],\n  [\"512\u00f72=256, 0\", \"802\u00f74=200, 2\"],\n  [\"424\u00f72=212, 0\", \"671\u00f78=83, 7\"],\n  [\"614\u00f74=153, 2\", \"150\u00f75=30, 0\"],\n  [\"538\u00f72=269, 0\", \"471\u00f72=235, 1\"],\n  [\"673\u00f79=74, 7\", \"450\u00f74=112, 2\"],\n  [\"756\u00f72=378, 0\", \"880\u00f75=176, 0\"],\n  [\"653\u00f76=108, 5\", \"634\u00f79=70, 4\"],\n  [\"289\u00f78=36, 1\", \"330\u00f79=36, 6\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text: ${oldText}`);\n  }\n\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n\n", "ps1": "# Update the three-digit-by-one-digit division answers in the table.\n# Each \"old\" answer string is unique within the document, so a plain\n# (non-whole-word, case-sensitive) Find/Replace of each pair is safe\n# and unambiguous. Old/new values are kept as two parallel arrays\n# (nested array literals inside @(...) are not reliable in this host).\n\n$d = $word.ActiveDocument\n\n$oldValues = @(\n    \"766\u00f78=95, 6\",\n    \"916\u00f76=152, 4\",\n    \"375\u00f74=93, 3\",\n    \"825\u00f75=165, 0\",\n    \"936\u00f76=156, 0\",\n    \"279\u00f77=39, 6\",\n    \"851\u00f76=141, 5\",\n    \"351\u00f74=87, 3\",\n    \"152\u00f73=50, 2\",\n    \"692\u00f76=115, 2\",\n    \"432\u00f75=86, 2\",\n    \"772\u00f72=386, 0\",\n    \"791\u00f72=395, 1\",\n    \"498\u00f74=124, 2\",\n    \"786\u00f76=131, 0\",\n    \"116\u00f79=12, 8\",\n    \"583\u00f79=64, 7\",\n    \"512\u00f72=256, 0\",\n    \"424\u00f72=212, 0\",\n    \"614\u00f74=153, 2\",\n    \"538\u00f72=269, 0\",\n    \"673\u00f79=74, 7\",\n    \"756\u00f72=378, 0\",\n    \"653\u00f76=108, 5\",\n    \"289\u00f78=36, 1\"\n)\n\n$newValues = @(\n    \"892\u00f75=178, 2\",\n    \"956\u00f73=318, 2\",\n    \"644\u00f72=322, 0\",\n    \"545\u00f76=90, 5\",\n    \"204\u00f74=51, 0\",\n    \"250\u00f75=50, 0\",\n    \"788\u00f78=98, 4\",\n    \"781\u00f77=111, 4\",\n    \"428\u00f73=142, 2\",\n    \"672\u00f74=168, 0\",\n    \"514\u00f79=57, 1\",\n    \"837\u00f77=119, 4\",\n    \"114\u00f78=14, 2\",\n    \"457\u00f72=228, 1\",\n    \"865\u00f74=216, 1\",\n    \"515\u00f76=85, 5\",\n    \"398\u00f73=132, 2\",\n    \"802\u00f74=200, 2\",\n    \"671\u00f78=83, 7\",\n    \"150\u00f75=30, 0\",\n    \"471\u00f72=235, 1\",\n    \"450\u00f74=112, 2\",\n    \"880\u00f75=176, 0\",\n    \"634\u00f79=70, 4\",\n    \"330\u00f79=36, 6\"\n)\n\nfor ($i = 0; $i -lt $oldValues.Count; $i++) {\n    $oldText = $oldValues[$i]\n    $newText = $newValues[$i]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $found = $find.Execute(\n        $oldText,   # FindText\n        $true,      # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $newText,   # ReplaceWith\n        2           # Replace (wdReplaceAll)\n    )\n\n    if (-not $found) {\n        throw (\"Could not find text at index \" + $i)\n    }\n}\n\n"}
